$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: finish filling in the existing date row (date changes from 7/20 to 7/25/2016)
$ws.Range("A3").Value = 42576
$ws.Range("B3").Value = 5.3
$ws.Range("C3").Value = 6
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 30
$ws.Range("F3").Value = "pruebas unitarias"

# Row 4: brand-new row (7/26/2016) - copy A3's date style down first so A4 keeps the
# same date number format (style index) instead of Excel minting a new one
$ws.Range("A3").Copy($ws.Range("A4"))
$ws.Range("A4").Value = 42577
$ws.Range("B4").Value = 6.3
$ws.Range("C4").Value = 7.1
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 20
$ws.Range("F4").Value = "profiler"

$ws.Range("C4").Select()
